$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59:103 down to 60:104
$ws.Rows.Item(59).EntireRow.Insert()

# Populate the newly inserted row 59 with the new weekly record
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value = 45233
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = 100112031
$ws.Cells.Item(59, 7).Value = "Poroto verde"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 705
$ws.Cells.Item(59, 11).Value = 1100
$ws.Cells.Item(59, 12).Value = 1100
$ws.Cells.Item(59, 13).Value = 1100
$ws.Cells.Item(59, 14).Value = "`$/kilo"
$ws.Cells.Item(59, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 16).Value = 1100
$ws.Cells.Item(59, 17).Value = 1
$ws.Cells.Item(59, 18).Value = "Hortaliza"
